$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '98.500.07'
$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("D3").Value = '3.359.12'
$ws.Range("E3").Value = '  +0.38%  '

Set-TextValue $ws.Range("D5") '256.27'
$ws.Range("E5").Value = '  -1.16%  '

Set-TextValue $ws.Range("D6") '663.53'
$ws.Range("E6").Value = '  +6.06%  '

$ws.Range("E7").Value = '  +5.60%  '

Set-TextValue $ws.Range("D8") '0.471'
$ws.Range("E8").Value = '  +19.52%  '

$ws.Range("E9").Value = '  +18.93%  '

$ws.Range("E10").Value = '  -0.03%  '

$ws.Range("D11").Value = '3.356.76'
$ws.Range("E11").Value = '  +0.38%  '

Set-TextValue $ws.Range("D12") '0.215'
$ws.Range("E12").Value = '  +7.93%  '

Set-TextValue $ws.Range("D13") '42.30'
$ws.Range("E13").Value = '  +11.11%  '

$ws.Range("E14").Value = '  +8.65%  '

$ws.Range("D15").Value = '98.623.68'
$ws.Range("E15").Value = '  -0.03%  '

$ws.Range("D16").Value = '3.992.80'
$ws.Range("E16").Value = '  +0.65%  '

Set-TextValue $ws.Range("D17") '5.67'
$ws.Range("E17").Value = '  +2.68%  '

Set-TextValue $ws.Range("D18") '7.95'
$ws.Range("E18").Value = '  +28.81%  '

$ws.Range("D19").Value = '3.352.06'
$ws.Range("E19").Value = '  +0.46%  '

Set-TextValue $ws.Range("D20") '17.01'
$ws.Range("E20").Value = '  +10.86%  '

Set-TextValue $ws.Range("D21") '528.11'
$ws.Range("E21").Value = '  +7.45%  '

Set-TextValue $ws.Range("D22") '3.55'
$ws.Range("E22").Value = '  -0.92%  '

Set-TextValue $ws.Range("D23") '10.47'
$ws.Range("E23").Value = '  +10.53%  '

$ws.Range("E24").Value = '  +1.78%  '

Set-TextValue $ws.Range("D25") '0.444'
$ws.Range("E25").Value = '  +48.53%  '

Set-TextValue $ws.Range("D26") '103.01'
$ws.Range("E26").Value = '  +14.53%  '

$ws.Range("E27").Value = '  +10.06%  '

Set-TextValue $ws.Range("D28") '12.60'
$ws.Range("E28").Value = '  +4.73%  '

$ws.Range("D29").Value = '3.536.28'
$ws.Range("E29").Value = '  +0.49%  '

$ws.Range("E30").Value = '  +9.24%  '

Set-TextValue $ws.Range("D31") '0.999'
$ws.Range("E31").Value = '  -0.12%  '

Set-TextValue $ws.Range("D32") '11.24'
$ws.Range("E32").Value = '  +15.07%  '

Set-TextValue $ws.Range("D33") '0.191'
$ws.Range("E33").Value = '  -0.67%  '

$ws.Range("E34").Value = '  -0.18%  '

Set-TextValue $ws.Range("D35") '29.51'
$ws.Range("E35").Value = '  +3.74%  '

Set-TextValue $ws.Range("D36") '0.545'
$ws.Range("E36").Value = '  +17.44%  '

$ws.Range("E37").Value = '  +8.74%  '

Set-TextValue $ws.Range("D38") '7.73'
$ws.Range("E38").Value = '  +5.82%  '

$ws.Range("E39").Value = '  +4.66%  '

Set-TextValue $ws.Range("D40") '521.99'
$ws.Range("E40").Value = '  +3.76%  '

Set-TextValue $ws.Range("D41") '3.94'
$ws.Range("E41").Value = '  +7.28%  '

$ws.Range("E42").Value = '  +5.83%  '

Set-TextValue $ws.Range("D43") '24.73'
$ws.Range("E43").Value = '  -0.63%  '

$ws.Range("E44").Value = '  +31.65%  '

Set-TextValue $ws.Range("D45") '0.830'
$ws.Range("E45").Value = '  +6.01%  '

$ws.Range("E46").Value = '  +1.28%  '

$ws.Range("E47").Value = '  -0.04%  '

Set-TextValue $ws.Range("D48") '5.19'
$ws.Range("E48").Value = '  +10.31%  '

$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D49") '1.55'
$ws.Range("E49").Value = '  +13.59%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D50") '2.07'
$ws.Range("E50").Value = '  +5.19%  '

Set-TextValue $ws.Range("D51") '7.91'
$ws.Range("E51").Value = '  +14.13%  '
